$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tickers = @(
    "AAF", "ABDN", "ABF", "ANTO", "AUTO", "AV", "BARC", "BATS", "BDEV", "BEZ",
    "BF.B", "BKG", "BNZL", "BRBY", "BRK.B", "BT-A", "CCH", "CRDA", "DCC", "DGE",
    "ENT", "EXPN", "FCIT", "FRAS", "FRES", "GLEN", "HLMA", "HSBA", "HSX", "IMB",
    "INF", "ITRK", "JMAT", "KGF", "LGEN", "LLOY", "LSEG", "MNDI", "MNG", "OCDO",
    "PHNX", "PSH", "PSON", "REL", "RMV", "RR", "RS1", "SBRY", "SDR", "SGE",
    "SGRO", "SKG", "SMDS", "SMT", "SN", "SPX", "SSE", "STAN", "STJ", "ULVR",
    "UU", "WEIR", "WTB"
)

$startRow = 5
for ($i = 0; $i -lt $tickers.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $tickers[$i]
}
